$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F206").Value = "root.pop_gridConnections[221]"
$ws.Range("F207").Value = "root.pop_gridConnections[222]"
$ws.Range("F208").Value = "root.pop_gridConnections[223]"
$ws.Range("F209").Value = "root.pop_gridConnections[224]"
$ws.Range("F210").Value = "root.pop_gridConnections[225]"
$ws.Range("F211").Value = "root.pop_gridConnections[226]"
$ws.Range("F212").Value = "root.pop_gridConnections[227]"
$ws.Range("F213").Value = "root.pop_gridConnections[228]"
$ws.Range("F214").Value = "root.pop_gridConnections[229]"
$ws.Range("F215").Value = "root.pop_gridConnections[230]"
$ws.Range("F216").Value = "root.pop_gridConnections[231]"
$ws.Range("F217").Value = "root.pop_gridConnections[232]"
$ws.Range("F218").Value = "root.pop_gridConnections[233]"
$ws.Range("F219").Value = "root.pop_gridConnections[234]"
$ws.Range("F220").Value = "root.pop_gridConnections[235]"
$ws.Range("F221").Value = "root.pop_gridConnections[236]"
$ws.Range("F222").Value = "root.pop_gridConnections[237]"
$ws.Range("F223").Value = "root.pop_gridConnections[238]"
$ws.Range("F224").Value = "root.pop_gridConnections[239]"
$ws.Range("F225").Value = "root.pop_gridConnections[240]"
$ws.Range("F226").Value = "root.pop_gridConnections[241]"
$ws.Range("F227").Value = "root.pop_gridConnections[242]"
$ws.Range("F228").Value = "root.pop_gridConnections[243]"
$ws.Range("F229").Value = "root.pop_gridConnections[244]"
$ws.Range("F230").Value = "root.pop_gridConnections[245]"
$ws.Range("F231").Value = "root.pop_gridConnections[246]"
$ws.Range("F232").Value = "root.pop_gridConnections[247]"
$ws.Range("F233").Value = "root.pop_gridConnections[248]"
$ws.Range("F234").Value = "root.pop_gridConnections[249]"
$ws.Range("F235").Value = "root.pop_gridConnections[250]"
$ws.Range("F236").Value = "root.pop_gridConnections[251]"
$ws.Range("F237").Value = "root.pop_gridConnections[252]"
$ws.Range("F238").Value = "root.pop_gridConnections[253]"
$ws.Range("F239").Value = "root.pop_gridConnections[254]"
$ws.Range("F240").Value = "root.pop_gridConnections[255]"
$ws.Range("F241").Value = "root.pop_gridConnections[256]"
$ws.Range("F242").Value = "root.pop_gridConnections[257]"
$ws.Range("F243").Value = "root.pop_gridConnections[258]"
$ws.Range("F244").Value = "root.pop_gridConnections[259]"
$ws.Range("F245").Value = "root.pop_gridConnections[260]"
$ws.Range("F246").Value = "root.pop_gridConnections[261]"
$ws.Range("F247").Value = "root.pop_gridConnections[262]"
$ws.Range("F248").Value = "root.pop_gridConnections[263]"
$ws.Range("F249").Value = "root.pop_gridConnections[264]"
$ws.Range("F250").Value = "root.pop_gridConnections[265]"
$ws.Range("F251").Value = "root.pop_gridConnections[266]"
$ws.Range("F252").Value = "root.pop_gridConnections[267]"
$ws.Range("F253").Value = "root.pop_gridConnections[268]"
$ws.Range("F254").Value = "root.pop_gridConnections[269]"
$ws.Range("F255").Value = "root.pop_gridConnections[270]"
$ws.Range("F256").Value = "root.pop_gridConnections[271]"
$ws.Range("F257").Value = "root.pop_gridConnections[272]"
$ws.Range("F258").Value = "root.pop_gridConnections[273]"
$ws.Range("F259").Value = "root.pop_gridConnections[274]"
$ws.Range("F260").Value = "root.pop_gridConnections[275]"
$ws.Range("F261").Value = "root.pop_gridConnections[276]"
$ws.Range("F262").Value = "root.pop_gridConnections[277]"
$ws.Range("F263").Value = "root.pop_gridConnections[278]"
$ws.Range("F264").Value = "root.pop_gridConnections[279]"
$ws.Range("F265").Value = "root.pop_gridConnections[280]"
$ws.Range("F266").Value = "root.pop_gridConnections[281]"
$ws.Range("F267").Value = "root.pop_gridConnections[282]"
$ws.Range("F268").Value = "root.pop_gridConnections[283]"
$ws.Range("F269").Value = "root.pop_gridConnections[284]"
$ws.Range("F270").Value = "root.pop_gridConnections[285]"
$ws.Range("F271").Value = "root.pop_gridConnections[286]"
$ws.Range("F272").Value = "root.pop_gridConnections[287]"
$ws.Range("F273").Value = "root.pop_gridConnections[288]"
$ws.Range("F274").Value = "root.pop_gridConnections[289]"
$ws.Range("F275").Value = "root.pop_gridConnections[290]"
$ws.Range("F276").Value = "root.pop_gridConnections[291]"
$ws.Range("F277").Value = "root.pop_gridConnections[292]"
$ws.Range("F278").Value = "root.pop_gridConnections[293]"
$ws.Range("F279").Value = "root.pop_gridConnections[294]"
$ws.Range("F280").Value = "root.pop_gridConnections[295]"
$ws.Range("F281").Value = "root.pop_gridConnections[296]"
$ws.Range("F282").Value = "root.pop_gridConnections[297]"
$ws.Range("F283").Value = "root.pop_gridConnections[298]"

$ws.Range("P206").Value = "'-342.75986356469923"
$ws.Range("P207").Value = "'-382.0718635647259"
$ws.Range("P208").Value = "'-187.69586356462963"
$ws.Range("P209").Value = "'-633.2318635648804"
$ws.Range("P210").Value = "'-770.151863564898"
$ws.Range("P211").Value = "'-1039.3718635649357"
$ws.Range("P212").Value = "'-1275.7478635649186"
$ws.Range("P213").Value = "'-439.19186356476405"
$ws.Range("P214").Value = "'-266.3198635646405"
$ws.Range("P215").Value = "'-277.23986356464917"
$ws.Range("P216").Value = "'-280.3058635646509"
$ws.Range("P217").Value = "'-379.8878635647259"
$ws.Range("P218").Value = "'-663.6398635648867"
$ws.Range("P219").Value = "'-168.0398635646302"
$ws.Range("P220").Value = "'-139.64786356463966"
$ws.Range("P221").Value = "'-287.73986356465736"
$ws.Range("P222").Value = "'-518.9918635648086"
$ws.Range("P223").Value = "'-395.1758635647387"
$ws.Range("P224").Value = "'-255.3998635646296"
$ws.Range("P225").Value = "'-880.0238635649142"
$ws.Range("P226").Value = "'-746.2538635648955"
$ws.Range("P227").Value = "'-135.27986356464095"
$ws.Range("P228").Value = "'-2929.8758635647346"
$ws.Range("P229").Value = "'-2013.939863564789"
$ws.Range("P230").Value = "'-268.7978635646423"
$ws.Range("P231").Value = "'-950.2478635649197"
$ws.Range("P232").Value = "'-353.6798635647065"
$ws.Range("P233").Value = "'-246.0758635646293"
$ws.Range("P234").Value = "'-413.1938635647472"
$ws.Range("P235").Value = "'-372.2438635647198"
$ws.Range("P236").Value = "'-327.26186356468696"
$ws.Range("P237").Value = "'-368.96786356471677"
$ws.Range("P238").Value = "'-414.4958635647474"
$ws.Range("P239").Value = "'-185.51186356462998"
$ws.Range("P240").Value = "'-417.477863564753"
$ws.Range("P241").Value = "'-2524.9958635647026"
$ws.Range("P242").Value = "'-209.5358635646299"
$ws.Range("P243").Value = "'-152.7518635646358"
$ws.Range("P244").Value = "'-234.14786356462955"
$ws.Range("P245").Value = "'-490.1798635647937"
$ws.Range("P246").Value = "'-202.98386356463016"
$ws.Range("P247").Value = "'-644.9078635648843"
$ws.Range("P248").Value = "'-401.7278635647412"
$ws.Range("P249").Value = "'-164.2178635646303"
$ws.Range("P250").Value = "'-547.7198635648267"
$ws.Range("P251").Value = "'-146.19986356463835"
$ws.Range("P252").Value = "'-272.87186356464616"
$ws.Range("P253").Value = "'-240.11186356462983"
$ws.Range("P254").Value = "'-237.92786356462935"
$ws.Range("P255").Value = "'-266.31986356463926"
$ws.Range("P256").Value = "'-403.91186356474157"
$ws.Range("P257").Value = "'-198.61586356462976"
$ws.Range("P258").Value = "'-711.8558635648923"
$ws.Range("P259").Value = "'-340.5758635646977"
$ws.Range("P260").Value = "'-545.4518635648261"
$ws.Range("P261").Value = "'-622.6478635648753"
$ws.Range("P262").Value = "'-824.6678635649091"
$ws.Range("P263").Value = "'-707.4458635648913"
$ws.Range("P264").Value = "'-320.9198635646831"
$ws.Range("P265").Value = "'-1344.7958635649081"
$ws.Range("P266").Value = "'-462.41786356477627"
$ws.Range("P267").Value = "'-667.7138635648847"
$ws.Range("P268").Value = "'-459.0578635647758"
$ws.Range("P269").Value = "'-130.91186356464146"
$ws.Range("P270").Value = "'-508.9118635648058"
$ws.Range("P271").Value = "'-695.6018635648871"
$ws.Range("P272").Value = "'-850.9598635649144"
$ws.Range("P273").Value = "'-1098.4238635649417"
$ws.Range("P274").Value = "'-1137.1898635649397"
$ws.Range("P275").Value = "'-137.46386356464026"
$ws.Range("P276").Value = "'-478.1678635647865"
$ws.Range("P277").Value = "'-452.54786356476984"
$ws.Range("P278").Value = "'-216.08786356462943"
$ws.Range("P279").Value = "'-269.133863564642"
$ws.Range("P280").Value = "'-817.1918635649076"
$ws.Range("P281").Value = "'-240.11186356462963"
$ws.Range("P282").Value = "'-209.87186356462982"
$ws.Range("P283").Value = "'-185.51186356462998"
